$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "ToC"
$ws.Range("B13").Value = "Table of Content"
$ws.Range("A14").Value = "HIES"
$ws.Range("B14").Value = "Household Income and Expenditure Survey"

$ws.Range("B15").Select()
